$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update File Name column (A2:A8)
$ws.Range("A2").Value = "ck20.txt"
$ws.Range("A3").Value = "ck15.txt"
$ws.Range("A4").Value = "ck10.txt"
$ws.Range("A5").Value = "ck00.txt"
$ws.Range("A6").Value = "ckm10.txt"
$ws.Range("A7").Value = "ckm15.txt"
$ws.Range("A8").Value = "ckm20.txt"

# Update Clok_value[V] column (B2:B8)
$ws.Range("B2").Value = 2
$ws.Range("B3").Value = 1.5
$ws.Range("B4").Value = 1
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = -1
$ws.Range("B7").Value = -1.5
$ws.Range("B8").Value = -2

# Update # Vin Values column (C2:C8)
$ws.Range("C2").Value = 24
$ws.Range("C3").Value = 24
$ws.Range("C4").Value = 24
$ws.Range("C5").Value = 24
$ws.Range("C6").Value = 24
$ws.Range("C7").Value = 24
$ws.Range("C8").Value = 24

$ws.Range("H18").Select()
